$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the one remaining data row's values. These are stored as plain text
# (shared strings) in the workbook, so a leading apostrophe forces the
# number-looking ones ("2022", "1234567890") to stay text instead of turning
# into numeric cells; resetting the style back to Normal afterwards drops the
# quote-prefix formatting so the cell keeps its original (default) style.
$ws.Range("A2").Value = "'2022"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").Value = "btcs"

$ws.Range("D2").Value = "'1234567890"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "Institute 1"

# Delete rows 3 through 6 (the extra dropdown sample rows), keeping header + first data row
$ws.Range("A3:E6").EntireRow.Delete()

$wb.Save()
